$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.224.72"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.199.25"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.60"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.87"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.91"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.88"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.84"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "2.523.07"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "2.193.43"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.793"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("D19").Value = "42.010.76"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.55"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.87"
$ws.Range("E22").Value = "  -5.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.07"
$ws.Range("E23").Value = "  -9.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "228.14"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.31"
$ws.Range("E28").Value = "  -8.59%  "
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.14"
$ws.Range("E31").Value = "  +2.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.07"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.07"
$ws.Range("E33").Value = "  +10.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0786"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  -4.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0317"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.52"
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.44"
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.20"
$ws.Range("E43").Value = "  -7.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.193"
$ws.Range("E44").Value = "  -3.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.50"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0973"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.84"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.414"
$ws.Range("E51").Value = "  +11.88%  "
